{"js": "// Replace the overfitting-ratio paragraph text. The paragraph has a\n// proofErr-wrapped run around the word \"far\" in the middle, so the\n// replacement is split into two search/replace operations that bracket\n// that run without touching it.\nconst seg1Search = context.document.body.search(\n  \"Whilst training I output an overfitting ratio: train loss / test loss. This gives a good idea of whether the model is going too \",\n  { matchCase: true }\n);\nseg1Search.load(\"items\");\nawait context.sync();\nif (seg1Search.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for paragraph 1 segment 1, found \" + seg1Search.items.length);\n}\nseg1Search.items[0].insertText(\n  \"Whilst training, an overfitting ratio, defined as test loss / train loss, is computed. For example, if the loss on the training set is half that on the test set, overfit will be equal to two. This gives a good indication of whether the model is going too \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nconst seg2Search = context.document.body.search(\n  \" and we need to stop training.\",\n  { matchCase: true }\n);\nseg2Search.load(\"items\");\nawait context.sync();\nif (seg2Search.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for paragraph 1 segment 2, found \" + seg2Search.items.length);\n}\nseg2Search.items[0].insertText(\n  \" and we need to stop training, if overfit > 1.3 for example.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// Second paragraph: \"it would be nice to be able to reproduce important\" -> \"it is desirable to reproduce\"\nconst para2Search = context.document.body.search(\n  \"However, for this audio generation problem, it would be nice to be able to reproduce important samples as accurately as possible, so overfitting the original training samples may not necessarily be a bad thing.\",\n  { matchCase: true }\n);\npara2Search.load(\"items\");\nawait context.sync();\nif (para2Search.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for paragraph 2, found \" + para2Search.items.length);\n}\npara2Search.items[0].insertText(\n  \"However, for this audio generation problem, it is desirable to reproduce samples as accurately as possible, so overfitting the original training samples may not necessarily be a bad thing.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// Third paragraph: \"generate diverse outputs\" -> \"generalise to samples\"\nconst para3Search = context.document.body.search(\n  \"Over-fitting may also lead to the Auto-Encoder being less able to generate diverse outputs outside the original training dataset. \",\n  { matchCase: true }\n);\npara3Search.load(\"items\");\nawait context.sync();\nif (para3Search.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for paragraph 3, found \" + para3Search.items.length);\n}\npara3Search.items[0].insertText(\n  \"Over-fitting may also lead to the Auto-Encoder being less able to generalise to samples outside the original training dataset. \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    $result = $find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2)\n    if (-not $result) {\n        throw \"Replace failed for: $searchText\"\n    }\n}\n\n# Paragraph 1 is split around a proofErr-wrapped run containing just the\n# word \"far\" - replace the text before and after it separately so that\n# run is left untouched.\nReplace-Text \"Whilst training I output an overfitting ratio: train loss / test loss. This gives a good idea of whether the model is going too \" \"Whilst training, an overfitting ratio, defined as test loss / train loss, is computed. For example, if the loss on the training set is half that on the test set, overfit will be equal to two. This gives a good indication of whether the model is going too \"\n\nReplace-Text \" and we need to stop training.\" \" and we need to stop training, if overfit > 1.3 for example.\"\n\n# Paragraph 2\nReplace-Text \"However, for this audio generation problem, it would be nice to be able to reproduce important samples as accurately as possible, so overfitting the original training samples may not necessarily be a bad thing.\" \"However, for this audio generation problem, it is desirable to reproduce samples as accurately as possible, so overfitting the original training samples may not necessarily be a bad thing.\"\n\n# Paragraph 3\nReplace-Text \"Over-fitting may also lead to the Auto-Encoder being less able to generate diverse outputs outside the original training dataset. \" \"Over-fitting may also lead to the Auto-Encoder being less able to generalise to samples outside the original training dataset. \"\n"}
